# Generate Report for Handback
# Updates the "Correspond Handback DateTime"-adjacent timestamp cells
# (Correspond Handoff Datetime / Correspond Handback DateTime) for the
# 5951cf53-b35b-4d3c-89a6-9acd91344a56 row on both the zh-cn and de-de
# status sheets to reflect a freshly (re-)generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-17 20:34:26"
$wsZhCn.Range("H4").Value = "2016-03-17 20:34:48"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-17 20:34:33"
$wsDeDe.Range("H4").Value = "2016-03-17 20:34:54"
